$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 26 (shifts rows 26-37 down to 27-38)
$ws.Rows.Item(26).Insert()

$ws.Cells.Item(26, 1).Value = "brd_id"
$ws.Cells.Item(26, 2).Value = 4
$ws.Cells.Item(26, 3).Value = "cosmic_ray_eth_control"
$ws.Cells.Item(26, 4).Value = 3
$ws.Cells.Item(26, 5).Value = 32
$ws.Cells.Item(26, 6).Value = "readwrite"
$ws.Cells.Item(26, 7).Value = "Number to identify each SNAP2 board"

# cr_dest_port (now row 27) shifted to start right after brd_id: offset_from_msb 3 -> 7
$ws.Cells.Item(27, 4).Value = 7

$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1
$ws.Range("D27").Select()
